$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1036.5625
$ws.Range("I19").Value = 765.8333
$ws.Range("J19").Value = 1199
$ws.Range("K19").Value = 765.8333
$ws.Range("L19").Value = 1199
$ws.Range("M19").Value = -590.8333
$ws.Range("N19").Value = -1549

$ws.Range("H53").Value = 172.95833
$ws.Range("J53").Value = 194.85715
$ws.Range("L53").Value = 194.85715
$ws.Range("N53").Value = -1468.85715

$ws.Range("H54").Value = 5000
$ws.Range("I54").Value = 5000
$ws.Range("K54").Value = 5000
$ws.Range("M54").Value = -4514

$ws.Range("H70").Value = 1603.2222
$ws.Range("I70").Value = 2700
$ws.Range("J70").Value = 1538.7059
$ws.Range("K70").Value = 8100
$ws.Range("L70").Value = 4616.1177
$ws.Range("M70").Value = -7830
$ws.Range("N70").Value = -5156.1177

$ws.Range("H73").Value = 1603.2222
$ws.Range("I73").Value = 2700
$ws.Range("J73").Value = 1538.7059
$ws.Range("K73").Value = 8100
$ws.Range("L73").Value = 4616.1177
$ws.Range("M73").Value = -7164
$ws.Range("N73").Value = -6488.1177

$ws.Range("H92").Value = 1000
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 1000
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -3496

$ws.Range("H108").Value = 37200
$ws.Range("J108").Value = 37200
$ws.Range("L108").Value = 37200
$ws.Range("N108").Value = -44880

$ws.Range("H116").Value = 2349.5293
$ws.Range("I116").Value = 2040.909
$ws.Range("J116").Value = 2915.3333
$ws.Range("K116").Value = 2040.909
$ws.Range("L116").Value = 2915.3333
$ws.Range("M116").Value = 1401.091
$ws.Range("N116").Value = -9799.3333

$ws.Range("H129").Value = 796.6061
$ws.Range("I129").Value = 295.2857
$ws.Range("J129").Value = 931.5769
$ws.Range("K129").Value = 885.8571000000001
$ws.Range("L129").Value = 2794.7307
$ws.Range("M129").Value = 4114.1429
$ws.Range("N129").Value = -12794.7307

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19315.418
$ws.Range("I32").Value = 20257.678
$ws.Range("K32").Value = 20257.678
$ws.Range("M32").Value = -19970.678

$ws.Range("H61").Value = 6932.3
$ws.Range("I61").Value = 3990.1428
$ws.Range("J61").Value = 20038.273
$ws.Range("K61").Value = 3990.1428
$ws.Range("L61").Value = 20038.273
$ws.Range("M61").Value = -3778.1428
$ws.Range("N61").Value = -20462.273

$ws.Range("H74").Value = 7163.4473
$ws.Range("I74").Value = 4976.5757
$ws.Range("J74").Value = 21596.8
$ws.Range("K74").Value = 4976.5757
$ws.Range("L74").Value = 21596.8
$ws.Range("M74").Value = -4102.5757
$ws.Range("N74").Value = -23344.8

$ws.Range("H77").Value = 7163.4473
$ws.Range("I77").Value = 4976.5757
$ws.Range("J77").Value = 21596.8
$ws.Range("K77").Value = 24882.8785
$ws.Range("L77").Value = 107984
$ws.Range("M77").Value = -20514.8785
$ws.Range("N77").Value = -116720

$ws.Range("H102").Value = 1236411.1
$ws.Range("I102").Value = 1544447.2
$ws.Range("J102").Value = 4266.5
$ws.Range("K102").Value = 1544447.2
$ws.Range("L102").Value = 4266.5
$ws.Range("M102").Value = -1542825.2
$ws.Range("N102").Value = -7510.5

$ws.Range("H132").Value = 5077.5527
$ws.Range("I132").Value = 1739.5
$ws.Range("J132").Value = 9667.375
$ws.Range("K132").Value = 5218.5
$ws.Range("L132").Value = 29002.125
$ws.Range("M132").Value = -2688.5
$ws.Range("N132").Value = -34062.125

$ws.Range("H136").Value = 6932.3
$ws.Range("I136").Value = 3990.1428
$ws.Range("J136").Value = 20038.273
$ws.Range("K136").Value = 11970.4284
$ws.Range("L136").Value = 60114.819
$ws.Range("M136").Value = -9420.428400000001
$ws.Range("N136").Value = -65214.819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 755.3333
$ws.Range("I20").Value = 803.8889
$ws.Range("J20").Value = 658.2222
$ws.Range("K20").Value = 803.8889
$ws.Range("L20").Value = 658.2222
$ws.Range("M20").Value = -556.8889
$ws.Range("N20").Value = -1152.2222

$ws.Range("H80").Value = 206.38461
$ws.Range("J80").Value = 209.60869
$ws.Range("L80").Value = 209.60869
$ws.Range("N80").Value = -2205.60869

$ws.Range("H83").Value = 206.38461
$ws.Range("J83").Value = 209.60869
$ws.Range("L83").Value = 1048.04345
$ws.Range("N83").Value = -11032.04345

$ws.Range("H86").Value = 2355.2632
$ws.Range("I86").Value = 2256.6667
$ws.Range("J86").Value = 2725
$ws.Range("K86").Value = 2256.6667
$ws.Range("L86").Value = 2725
$ws.Range("M86").Value = -1133.6667
$ws.Range("N86").Value = -4971

$ws.Range("H89").Value = 2355.2632
$ws.Range("I89").Value = 2256.6667
$ws.Range("J89").Value = 2725
$ws.Range("K89").Value = 11283.3335
$ws.Range("L89").Value = 13625
$ws.Range("M89").Value = -5667.333500000001
$ws.Range("N89").Value = -24857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2227.7546
$ws.Range("I31").Value = 1662.1052
$ws.Range("J31").Value = 3660.7334
$ws.Range("K31").Value = 1662.1052
$ws.Range("L31").Value = 3660.7334
$ws.Range("M31").Value = -1367.1052
$ws.Range("N31").Value = -4250.7334

$ws.Range("H34").Value = 2227.7546
$ws.Range("I34").Value = 1662.1052
$ws.Range("J34").Value = 3660.7334
$ws.Range("K34").Value = 1662.1052
$ws.Range("L34").Value = 3660.7334
$ws.Range("M34").Value = -1460.1052
$ws.Range("N34").Value = -4064.7334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5799.907
$ws.Range("I70").Value = 5299.5
$ws.Range("J70").Value = 6338.8076
$ws.Range("K70").Value = 5299.5
$ws.Range("L70").Value = 6338.8076
$ws.Range("M70").Value = -5029.5
$ws.Range("N70").Value = -6878.8076

$ws.Range("H73").Value = 5799.907
$ws.Range("I73").Value = 5299.5
$ws.Range("J73").Value = 6338.8076
$ws.Range("K73").Value = 5299.5
$ws.Range("L73").Value = 6338.8076
$ws.Range("M73").Value = -4363.5
$ws.Range("N73").Value = -8210.8076

$ws.Range("H132").Value = 5056.5713
$ws.Range("I132").Value = 2280.1875
$ws.Range("J132").Value = 34671.332
$ws.Range("K132").Value = 6840.5625
$ws.Range("L132").Value = 104013.996
$ws.Range("M132").Value = -4310.5625
$ws.Range("N132").Value = -109073.996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H54").Value = 6000
$ws.Range("I54").Value = 6000
$ws.Range("K54").Value = 6000
$ws.Range("M54").Value = -5356

$ws.Range("H130").Value = 290750
$ws.Range("J130").Value = 290750
$ws.Range("L130").Value = 290750
$ws.Range("N130").Value = -300790

$ws.Range("H132").Value = 3649.4834
$ws.Range("I132").Value = 3728.152
$ws.Range("J132").Value = 3391
$ws.Range("K132").Value = 11184.456
$ws.Range("L132").Value = 10173
$ws.Range("M132").Value = -8654.456
$ws.Range("N132").Value = -15233

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1014.4074
$ws.Range("I126").Value = 907.875
$ws.Range("J126").Value = 1866.6666
$ws.Range("K126").Value = 2723.625
$ws.Range("L126").Value = 5599.9998
$ws.Range("M126").Value = -253.625
$ws.Range("N126").Value = -10539.9998

$ws.Range("H132").Value = 1820.8889
$ws.Range("I132").Value = 883.9545000000001
$ws.Range("J132").Value = 3293.2144
$ws.Range("K132").Value = 2651.8635
$ws.Range("L132").Value = 9879.643199999999
$ws.Range("M132").Value = -121.8635000000004
$ws.Range("N132").Value = -14939.6432

$ws.Range("H136").Value = 3056.6235
$ws.Range("I136").Value = 1414.9246
$ws.Range("J136").Value = 5775.6875
$ws.Range("K136").Value = 4244.7738
$ws.Range("L136").Value = 17327.0625
$ws.Range("M136").Value = -1694.7738
$ws.Range("N136").Value = -22427.0625
